$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix doctor list: every "Chris Paul" entry at St Johns Hospital was actually
# a different doctor, "Toby Mac". Update the Name column (A) for each of
# those rows.
$rowsToFix = @(6, 15, 24, 33, 42, 51, 60, 69, 78, 87)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 1).Value = "Toby Mac"
}

# Give row 58 a slightly taller custom row height.
$ws.Rows.Item(58).RowHeight = 18

# Update the active selection to reflect the last edited cell.
$ws.Range("A87").Select()
